# Update gh-pages to output generated at 456a3b4
# Refreshes scraped numeric data ("想去人数" / "最低票价") across the
# four worksheets of the workbook.

$wb = $excel.ActiveWorkbook

function Set-Cells {
    param(
        [string]$SheetName,
        [hashtable]$Updates
    )
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($cellRef in $Updates.Keys) {
        $ws.Range($cellRef).Value = $Updates[$cellRef]
    }
}

Set-Cells "展览" @{
    "F3"  = 157
    "F4"  = 1788
    "G5"  = 106
    "F6"  = 1109
    "F7"  = 2237
    "F8"  = 2151
    "F9"  = 1123
    "F10" = 613
    "F11" = 24
    "F12" = 1688
    "G12" = 80
    "F13" = 405
    "F15" = 45
    "F17" = 235
    "F18" = 1603
    "F19" = 95
    "F20" = 653
    "F21" = 744
    "F22" = 96
    "F23" = 623
    "F24" = 12336
    "F25" = 12384
    "F26" = 918
    "F27" = 711
    "F30" = 30
    "F31" = 393
    "F32" = 1932
    "F35" = 208
    "F36" = 613
}

Set-Cells "演出" @{
    "F7" = 46
}

Set-Cells "本地生活" @{
    "F3" = 71
}

Set-Cells "全部类型" @{
    "F4"  = 157
    "F5"  = 1788
    "G6"  = 106
    "F7"  = 1109
    "F8"  = 2237
    "F9"  = 2151
    "F10" = 1123
    "F11" = 613
    "F12" = 71
    "F13" = 24
    "F14" = 1688
    "G14" = 80
    "F15" = 405
    "F18" = 45
    "F22" = 235
    "F23" = 1603
    "F24" = 95
    "F25" = 653
    "F26" = 744
    "F27" = 96
    "F28" = 623
    "F29" = 12336
    "F30" = 12384
    "F31" = 918
    "F32" = 711
    "F35" = 30
    "F36" = 393
    "F37" = 1933
    "F42" = 208
    "F43" = 613
    "F44" = 46
}

$wb.Save()
